# Additional companies sent for questionaire
# Remove the "Parent company" and "Location County/City" columns from the
# locomotive list, shifting remaining columns left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Parent company"); all following columns shift left by one.
$ws.Columns("B").Delete()

# After the above delete, the former column E ("Location County/City") is
# now column D. Delete it too, shifting the remaining columns left again.
$ws.Columns("D").Delete()

# Match the saved selection state from the edited workbook.
$ws.Range("A2:I2").Select()
